# Additions to White Paper, bibliography
#
# - C16: "DisplayPort, IEEE 1394" -> "DisplayPort, IEEE 1394, Decoders"
# - A1:  "topic" -> "domain" (old "topic" shared string becomes unused and
#        is dropped; "domain" is appended as a new shared string)
# - Selection moves from G31 to C27
#
# NOTE: order matters for shared-string layout parity with the source
# workbook - update C16 first (so the edited "DisplayPort..." string keeps
# its original slot) and only then touch A1 (whose old value is replaced,
# freeing its slot, with "domain" appended at the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "DisplayPort, IEEE 1394, Decoders"
$ws.Range("A1").Value = "domain"

[void]$ws.Range("C27").Select()
